# Add a new worksheet "CUMPLIMIENTO MENSUAL" after "VENTA MENSUAL", populate it
# with a per-group compliance summary for HIDALGO HIDALGO PEDRO GUSTAVO, and
# restore the originally-active sheet/tab afterwards.

$wb = $excel.ActiveWorkbook

$asesor = "HIDALGO HIDALGO PEDRO GUSTAVO"

$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)

# --- Create the new sheet as the 3rd tab (after "VENTA MENSUAL") ---------
$ws = $wb.Worksheets.Add($null, $sheet2)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# --- Header row ------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "ASESOR"
$ws.Cells.Item(1,2).Value = "GRUPO"
$ws.Cells.Item(1,3).Value = "PRESUPUESTO"
$ws.Cells.Item(1,4).Value = "VENTA"
$ws.Cells.Item(1,5).Value = "POR CUMPLIR"
$ws.Cells.Item(1,6).Value = "CUMPLIMIENTO"

# Match the bold/centered/bordered header style already used on the other
# sheets by copying the formatting from VENTAS POR GRUPO!A1.
$sheet1.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# --- Per-group data rows -----------------------------------------------------
$ws.Cells.Item(2,1).Value = $asesor
$ws.Cells.Item(2,2).Value = "240X120 PORCELANATO"
$ws.Cells.Item(2,3).Value = 782.465010521559
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 782.465010521559
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(3,1).Value = $asesor
$ws.Cells.Item(3,2).Value = "240X80 PORCELANATO"
$ws.Cells.Item(3,3).Value = 4168.07156573679
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 4168.07156573679
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(4,1).Value = $asesor
$ws.Cells.Item(4,2).Value = "FREGADEROS DE COCINA"
$ws.Cells.Item(4,3).Value = 513.831046659336
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 513.831046659336
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(5,1).Value = $asesor
$ws.Cells.Item(5,2).Value = "GRANITO"
$ws.Cells.Item(5,3).Value = 238.32
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 238.32
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(6,1).Value = $asesor
$ws.Cells.Item(6,2).Value = "GRIFERIAS"
$ws.Cells.Item(6,3).Value = 106.82
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 106.82
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(7,1).Value = $asesor
$ws.Cells.Item(7,2).Value = "INODOROS"
$ws.Cells.Item(7,3).Value = 1800
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 1800
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(8,1).Value = $asesor
$ws.Cells.Item(8,2).Value = "LAVABOS"
$ws.Cells.Item(8,3).Value = 625
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 625
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(9,1).Value = $asesor
$ws.Cells.Item(9,2).Value = "LED"
$ws.Cells.Item(9,3).Value = 300
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 300
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(10,1).Value = $asesor
$ws.Cells.Item(10,2).Value = "NO RESURTIBLES"
$ws.Cells.Item(10,3).Value = 650.25
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 650.25
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(11,1).Value = $asesor
$ws.Cells.Item(11,2).Value = "OTROS"
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 0
$ws.Cells.Item(11,5).Value = 0
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(12,1).Value = $asesor
$ws.Cells.Item(12,2).Value = "PANELES DECORATIVOS"
$ws.Cells.Item(12,3).Value = 350
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 350
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(13,1).Value = $asesor
$ws.Cells.Item(13,2).Value = "PANELES PU"
$ws.Cells.Item(13,3).Value = 230
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 230
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(14,1).Value = $asesor
$ws.Cells.Item(14,2).Value = "PANELES PVC"
$ws.Cells.Item(14,3).Value = 483
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 483
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(15,1).Value = $asesor
$ws.Cells.Item(15,2).Value = "PIEDRA SINTERIZADA"
$ws.Cells.Item(15,3).Value = 7465
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 7465
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(16,1).Value = $asesor
$ws.Cells.Item(16,2).Value = "PORCELANATO"
$ws.Cells.Item(16,3).Value = 29532.44
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 29532.44
$ws.Cells.Item(16,6).Value = 0
$ws.Cells.Item(17,1).Value = $asesor
$ws.Cells.Item(17,2).Value = "PUERTAS DE SEGURIDAD"
$ws.Cells.Item(17,3).Value = 342
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 342
$ws.Cells.Item(17,6).Value = 0
$ws.Cells.Item(18,1).Value = $asesor
$ws.Cells.Item(18,2).Value = "SAL SOLUBLE"
$ws.Cells.Item(18,3).Value = 2800
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 2800
$ws.Cells.Item(18,6).Value = 0

# --- TOTAL row ---------------------------------------------------------------
$ws.Cells.Item(19,2).Value = "TOTAL"
$ws.Cells.Item(19,3).Formula = "=SUM(C2:C18)"
$ws.Cells.Item(19,4).Formula = "=SUM(D2:D18)"
$ws.Cells.Item(19,5).Formula = "=SUM(E2:E18)"
$ws.Cells.Item(19,6).Value = 0

# --- Number formats: PRESUPUESTO / VENTA / POR CUMPLIR are currency, ---------
# --- CUMPLIMIENTO is a percentage (reuses the workbook's existing formats) --
$ws.Range("C2:E19").NumberFormat = """$""#,##0.00"
$ws.Range("F2:F19").NumberFormat = "0.00%"

# Right-align the "TOTAL" label.
$ws.Range("B19").HorizontalAlignment = -4152

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.166666666666668
$ws.Columns.Item(2).ColumnWidth = 21.166666666666668
$ws.Columns.Item(3).ColumnWidth = 21.166666666666668
$ws.Columns.Item(4).ColumnWidth = 10.166666666666668
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668

# --- Page margins (match the rest of the workbook's classic defaults) -------
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# --- Keep "VENTAS POR GRUPO" as the active sheet, as in the original file ---
$sheet1.Activate()
